$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column C (shifts the existing n_bio_reps/eps_mean/eps_mean_se
# columns one to the right, C/D/E -> D/E/F) to make room for the new
# "sigma_analytical" (analytical uncertainty) field between pN2 and n_bio_reps.
$ws.Range("C1").EntireColumn.Insert()

# Header for the new column.
$ws.Range("C1").Value = "sigma_analytical"

# Constant analytical-uncertainty value for every data row.
$ws.Range("C2:C9").Value = 0.084195259341794
